$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 316; this shifts the existing rows 316-357
# down to 317-358 and carries formatting (e.g. the date style on column D).
$ws.Rows("316:316").Insert()

# The columns that are constant for every record on this sheet (Mercado ID,
# Mercado, Region, Codreg, Tipo, Producto ID, Producto, Categoria ID,
# Categoria, Variedad) are copied from the row immediately below (row 317,
# which is the old row 316, pushed down by the insert above).
# NB: ".Value" misbehaves as a getter in this runtime, so read with
# ".Value2" instead.
$a = $ws.Cells.Item(317, 1).Value2
$b = $ws.Cells.Item(317, 2).Value2
$c = $ws.Cells.Item(317, 3).Value2
$e = $ws.Cells.Item(317, 5).Value2
$f = $ws.Cells.Item(317, 6).Value2
$g = $ws.Cells.Item(317, 7).Value2
$h = $ws.Cells.Item(317, 8).Value2
$i = $ws.Cells.Item(317, 9).Value2
$j = $ws.Cells.Item(317, 10).Value2
$k = $ws.Cells.Item(317, 11).Value2

# Populate the newly inserted row 316 with its data.
$ws.Range("A316").Value = $a
$ws.Range("B316").Value = $b
$ws.Range("C316").Value = $c
$ws.Range("D316").Value = 45142
$ws.Range("E316").Value = $e
$ws.Range("F316").Value = $f
$ws.Range("G316").Value = $g
$ws.Range("H316").Value = $h
$ws.Range("I316").Value = $i
$ws.Range("J316").Value = $j
$ws.Range("K316").Value = $k
$ws.Range("L316").Value = "Primera"
$ws.Range("M316").Value = 80
$ws.Range("N316").Value = 14000
$ws.Range("O316").Value = 14000
$ws.Range("P316").Value = 14000
$ws.Range("Q316").Value = "$/bandeja 18 kilos"
$ws.Range("R316").Value = "Región de O'Higgins"
$ws.Range("S316").Value = 778
$ws.Range("T316").Value = 18
